$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the per-task "Initial Estimate" values that shifted for this sprint.
$ws.Range("C5").Value = 4
$ws.Range("C7").Value = 3
$ws.Range("C13").Value = 2

# Row 17 (the "Fix testing problems" task) - estimate, week-1 remaining, and
# week-2 remaining all moved from 3 to 2.
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 2

# Recalculate dependent totals (SUM/SUMIF formulas in row 25 and G26:G28).
$excel.Calculate() | Out-Null

# Move the active selection to G28, matching where the editor ended up.
$ws.Range("G28").Select() | Out-Null
